$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, pushing the existing rows 112-176
# down to 113-177 (all their content, including styles, moves with them).
$ws.Rows(112).Insert()

# Populate the newly inserted row 112 with a new data record: a copy of
# the (now shifted) row 113 values, except the date (column D) and the
# volume (column J), which differ for this new entry.
$ws.Range("A112").Value = 5
$ws.Range("B112").Value = "Macroferia Regional de Talca"
$ws.Range("C112").Value = "Maule"
$ws.Range("D112").Value = 44488
$ws.Range("E112").Value = 7
$ws.Range("F112").Value = 100112003
$ws.Range("G112").Value = "Ajo"
$ws.Range("H112").Value = "Chino"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 300
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 15000
$ws.Range("M112").Value = 15000
$ws.Range("N112").Value = "$/malla 10 kilos"
$ws.Range("O112").Value = "China"
$ws.Range("P112").Value = 1500
$ws.Range("Q112").Value = 10
$ws.Range("R112").Value = "Hortaliza"
